$wb = $excel.ActiveWorkbook

# --- validLoginData sheet: update the login credentials used by the test data ---
$ws1 = $wb.Worksheets.Item("validLoginData")
$ws1.Range("A2").Value = "Admin"
$ws1.Range("B2").Value = "admin123"

# Restore the display text on the existing hyperlink (still pointing at the old
# mailto: target) now that the cell text itself has changed.
foreach ($hl in $ws1.Hyperlinks) {
    $hl.TextToDisplay = "Attitude665665@"
}

# This sheet becomes the active / selected tab, with a new selected cell.
$ws1.Activate() | Out-Null
$ws1.Range("D7").Select() | Out-Null
